$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the run that holds the inline picture (the triangle-of-stars image,
#    in the paragraph right after "Realizati urmatorul triunghi de stelute
#    ...") as NoProofing -> adds <w:noProof/> to that run's rPr.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.NoProofing = 1
        $found = $true
    }
}

# ---------------------------------------------------------------------------
# 2) Collapse the two empty "Listparagraf" placeholder paragraphs that sit
#    right after the "Cifrele unui numar" bullet item into a single plain
#    paragraph:
#      - the first empty bullet paragraph (style Listparagraf, no numbering)
#        is removed entirely;
#      - the second one (style Listparagraf + numPr ilvl0/numId1) loses its
#        pStyle/numPr but keeps its spacing/jc/rPr, i.e. becomes a bare
#        paragraph like the ones that already follow it.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()
    $style = $para.Range.ParagraphStyle.NameLocal
    $listType = $para.Range.ListFormat.ListType
    $shapes = $para.Range.InlineShapes.Count

    if ($text -eq "" -and $style -eq "List Paragraph" -and $listType -eq 0 -and $shapes -eq 0) {
        $next = $d.Paragraphs.Item($i + 1)
        $nextText = $next.Range.Text.Trim()
        $nextStyle = $next.Range.ParagraphStyle.NameLocal
        $nextListType = $next.Range.ListFormat.ListType

        if ($nextText -eq "" -and $nextStyle -eq "List Paragraph" -and $nextListType -eq 3) {
            # Delete the first (plain Listparagraf, no numbering) empty paragraph.
            $para.Range.Delete()

            # The following paragraph has shifted up to index $i; rewrite its
            # pPr so only spacing/jc/rPr remain (pStyle + numPr dropped).
            $target = $d.Paragraphs.Item($i).Range
            $target.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:after='0'/><w:jc w:val='both'/><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/><w:lang w:val='ro-RO'/></w:rPr></w:pPr></w:p>")
            break
        }
    }
}
